$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 319.35715
$ws.Range("I9").Value = 254.27272
$ws.Range("J9").Value = 558
$ws.Range("K9").Value = 254.27272
$ws.Range("L9").Value = 558
$ws.Range("M9").Value = -85.27271999999999
$ws.Range("N9").Value = -896

# ALC row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 801
$ws.Range("I18").Value = 801
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 801
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -517

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4752

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -4142

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5999
$ws.Range("I112").Value = 5000
$ws.Range("J112").Value = 6332
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 18996
$ws.Range("M112").Value = -13892
$ws.Range("N112").Value = -21212

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1610.2858
$ws.Range("I132").Value = 1518.9231
$ws.Range("J132").Value = 2798
$ws.Range("K132").Value = 4556.7693
$ws.Range("L132").Value = 8394
$ws.Range("M132").Value = -2026.7693

# ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1950
$ws.Range("I137").Value = 1950
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5850
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3300

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 676
$ws.Range("I2").Value = 676
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 676
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -563
$ws.Range("N2").ClearContents()

# ARM row 22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1833.3334
$ws.Range("I22").Value = 1833.3334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1833.3334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1534.3334

# ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 9322
$ws.Range("I41").Value = 651.6667
$ws.Range("J41").Value = 15824.75
$ws.Range("K41").Value = 651.6667
$ws.Range("L41").Value = 15824.75
$ws.Range("M41").Value = -237.6667
$ws.Range("N41").Value = -16652.75

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3595.125
$ws.Range("I45").Value = 1537.2858
$ws.Range("J45").Value = 18000
$ws.Range("K45").Value = 1537.2858
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = -1160.2858

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3274.75
$ws.Range("I61").Value = 2749.5
$ws.Range("J61").Value = 3800
$ws.Range("K61").Value = 2749.5
$ws.Range("L61").Value = 3800
$ws.Range("M61").Value = -2537.5
$ws.Range("N61").Value = -4224

# ARM row 86
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# ARM row 89
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 676
$ws.Range("I116").Value = 676
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 676
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1618
$ws.Range("N116").ClearContents()

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3274.75
$ws.Range("I136").Value = 2749.5
$ws.Range("J136").Value = 3800
$ws.Range("K136").Value = 8248.5
$ws.Range("L136").Value = 11400
$ws.Range("M136").Value = -5698.5
$ws.Range("N136").Value = -16500

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 676
$ws.Range("I3").Value = 676
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 676
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -562
$ws.Range("N3").ClearContents()

# BSM row 46
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 29449.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 29449.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 29449.5
$ws.Range("N46").Value = -30045.5

# CRP row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2381.2856
$ws.Range("I2").Value = 1754.6666
$ws.Range("J2").Value = 2851.25
$ws.Range("K2").Value = 1754.6666
$ws.Range("L2").Value = 2851.25
$ws.Range("M2").Value = -1641.6666
$ws.Range("N2").Value = -3077.25

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# CRP row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 38399.4
$ws.Range("I69").Value = 18000
$ws.Range("J69").Value = 51999
$ws.Range("K69").Value = 18000
$ws.Range("L69").Value = 51999
$ws.Range("M69").Value = -17251
$ws.Range("N69").Value = -53497

# CRP row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 38399.4
$ws.Range("I72").Value = 18000
$ws.Range("J72").Value = 51999
$ws.Range("K72").Value = 54000
$ws.Range("L72").Value = 155997
$ws.Range("M72").Value = -50256
$ws.Range("N72").Value = -163485

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 53998.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 53998.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 53998.5
$ws.Range("N74").Value = -55746.5

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 53998.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 53998.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 161995.5
$ws.Range("N77").Value = -170731.5

# CRP row 93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 28266.334
$ws.Range("I93").Value = 19999.5
$ws.Range("J93").Value = 44800
$ws.Range("K93").Value = 19999.5
$ws.Range("L93").Value = 44800
$ws.Range("M93").Value = -18127.5
$ws.Range("N93").Value = -48544

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 140.42857
$ws.Range("I11").Value = 182.8
$ws.Range("J11").Value = 34.5
$ws.Range("K11").Value = 548.4000000000001
$ws.Range("L11").Value = 103.5
$ws.Range("M11").Value = -408.4000000000001
$ws.Range("N11").Value = -383.5

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 107.818184
$ws.Range("I23").Value = 40.333332
$ws.Range("J23").Value = 133.125
$ws.Range("K23").Value = 120.999996
$ws.Range("L23").Value = 399.375
$ws.Range("M23").Value = 114.000004
$ws.Range("N23").Value = -869.375

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1627.75
$ws.Range("I80").Value = 1502.5
$ws.Range("J80").Value = 1753
$ws.Range("K80").Value = 1502.5
$ws.Range("L80").Value = 1753
$ws.Range("M80").Value = -504.5
$ws.Range("N80").Value = -3749

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1627.75
$ws.Range("I83").Value = 1502.5
$ws.Range("J83").Value = 1753
$ws.Range("K83").Value = 7512.5
$ws.Range("L83").Value = 8765
$ws.Range("M83").Value = -2520.5
$ws.Range("N83").Value = -18749

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1004
$ws.Range("I102").Value = 1004
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1004
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 618
$ws.Range("N102").ClearContents()

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9272.143
$ws.Range("I7").Value = 8811
$ws.Range("J7").Value = 10425
$ws.Range("K7").Value = 8811
$ws.Range("L7").Value = 10425
$ws.Range("M7").Value = -8699

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1983.1666
$ws.Range("I22").Value = 1749.6666
$ws.Range("J22").Value = 2216.6667
$ws.Range("K22").Value = 1749.6666
$ws.Range("L22").Value = 2216.6667
$ws.Range("M22").Value = -1454.6666
$ws.Range("N22").Value = -2806.6667

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1983.1666
$ws.Range("I27").Value = 1749.6666
$ws.Range("J27").Value = 2216.6667
$ws.Range("K27").Value = 1749.6666
$ws.Range("L27").Value = 2216.6667
$ws.Range("M27").Value = -1642.6666
$ws.Range("N27").Value = -2430.6667

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6165.2
$ws.Range("I32").Value = 1456.5
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 1456.5
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -1139.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9998
$ws.Range("I46").Value = 9998
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 9998
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -9810

# LTW row 53
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 11500
$ws.Range("I53").Value = 10000
$ws.Range("J53").Value = 13000
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 13000
$ws.Range("M53").Value = -9482
$ws.Range("N53").Value = -14036

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 9272.143
$ws.Range("I126").Value = 8811
$ws.Range("J126").Value = 10425
$ws.Range("K126").Value = 26433
$ws.Range("L126").Value = 31275
$ws.Range("M126").Value = -23963

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 724257.9
$ws.Range("I136").Value = 724257.9
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2172773.7
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2170223.7

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 50
